# Lab10 report revision
# - Reworked the "chart shows..." paragraph wording (Task 3 discussion)
# - Fixed a "then" -> "than" typo ("...more levels than if they were added...")
# - Split the trailing "Kyle O'Connor..." paragraph so the _GoBack bookmark
#   lives on its own paragraph, and dropped the now-duplicate empty paragraph.

$d = $word.ActiveDocument

# wdFindContinue = 1 ; wdReplaceAll = 2
$wdFindContinue = 1
$wdReplaceAll = 2

function Replace-Text($oldText, $newText) {
    $range = $d.Content
    $ok = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, `
                               $true, $wdFindContinue, $false, $newText, $wdReplaceAll)
    if (-not $ok) {
        throw "Find/Replace failed for: $oldText"
    }
}

# --- Paragraph: "The chart shows how many pointers ..." -------------------

Replace-Text "The chart shows how many pointers were followed to either find the node with" `
             "     The chart shows how many pointers were followed to either find the node with"

Replace-Text " the key or add a key to a node when adding 50 keys and then finding the same 50 keys." `
             " a key or add a key to a node. The code adds 50 keys and then finds the same 50 keys."

Replace-Text "Clearly finding the key looks more intensive for following pointers than adding does." `
             "From the table, clearly finding a key requires the code to follow more pointers than while adding."

Replace-Text "Now as for degree 3 through degree 5. The consistency" `
             "Now as for degree 3 through degree 5, the consistency"

# --- Paragraph: "I think the results of task 3 ..." : then -> than --------

Replace-Text "there will be more levels then if they were added to a tree of degree 5" `
             "there will be more levels than if they were added to a tree of degree 5"

# --- Split the "Kyle O'Connor ..." paragraph so the _GoBack bookmark moves -
#     onto its own paragraph, and drop the old trailing empty paragraph.

$apostrophe = [char]0x2019
$kyleRange = $d.Content
$kyleRange.Find.Execute("Kyle O" + $apostrophe + "Connor is responsible for the Lab Report") | Out-Null

# Delete the paragraph mark right after "...Lab Report" (and right after the
# zero-width bookmark that immediately follows it). That merges the old
# trailing empty paragraph into this one, so it disappears.
$endOfText = $kyleRange.End
$markRange = $d.Range($endOfText, $endOfText + 1)
$markRange.Delete()

# Now split right after "...Lab Report" so the bookmark (which stayed glued
# to the end of the paragraph) ends up alone on the new paragraph.
Replace-Text "Lab Report" "Lab Report^p"
